# Apply w:color (val=70AD47, themeColor=accent6) to both the paragraph-mark run
# properties (pPr/rPr) and the text run properties (r/rPr) of 8 specific list
# paragraphs, matching the target diff exactly. We rebuild each paragraph's full
# OOXML (preserving every existing attribute/element) with the <w:color> element
# inserted, then push it back in with Range.InsertXML so nothing else changes.

$d = $word.ActiveDocument

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="48A37575" w14:textId="5A6F736A" w:rsidR="00BA42E6" w:rsidRPr="009F7C6F" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>Implement a “Deal Cards” action, triggered by either a keyboard key or a UI button.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(13)
$p.Range.InsertXML($xml)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="41DDE165" w14:textId="780D7ACD" w:rsidR="00BA42E6" w:rsidRPr="009F7C6F" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>When triggered:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(14)
$p.Range.InsertXML($xml)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="781CBC00" w14:textId="4E92B1F2" w:rsidR="00BA42E6" w:rsidRPr="009F7C6F" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>The dealer plays an animation (e.g., reaching toward the table).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(15)
$p.Range.InsertXML($xml)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="016543D1" w14:textId="11C43D99" w:rsidR="00BA42E6" w:rsidRPr="009F7C6F" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>A sound effect plays (e.g., card dealing SFX).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(16)
$p.Range.InsertXML($xml)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="24E74720" w14:textId="55ED1A3F" w:rsidR="00BA42E6" w:rsidRPr="009F7C6F" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="16"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>A basic visual effect simulates card distribution to the player spots (sprite or 3D object).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(17)
$p.Range.InsertXML($xml)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3C6FB74C" w14:textId="1F4F01CD" w:rsidR="00BA42E6" w:rsidRPr="009F7C6F" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>Include a simple UI button labeled “Deal Cards.”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(19)
$p.Range.InsertXML($xml)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="24D7C437" w14:textId="2C48E531" w:rsidR="00BA42E6" w:rsidRPr="009F7C6F" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>Sync the animation and SFX when the button is pressed.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(20)
$p.Range.InsertXML($xml)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="17402E06" w14:textId="6F82AA8A" w:rsidR="00BA42E6" w:rsidRPr="00C63008" w:rsidRDefault="00BA42E6" w:rsidP="00C63008"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:before="240" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r w:rsidRPr="009F7C6F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>Implement round-based dealing: Deal one card at a time to each player.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p = $d.Paragraphs.Item(27)
$p.Range.InsertXML($xml)
